$p = $ppt.ActivePresentation

# The deck currently carries two themes:
#   theme1.xml -> "Office Theme" colours, used only by the Notes Master
#   theme2.xml -> "Integral" colours, used by the Slide Master / main design
#
# The target edit swaps which colour set backs the main design: the
# Slide Master's theme should take on the plain "Office Theme" palette
# (what used to live in theme1.xml) while the Integral palette moves off
# of the main design.
#
# The only theme whose colours are reachable/editable through the
# PowerPoint object model here is the one driving the active design
# (Slide.ThemeColorScheme / CustomLayout.ThemeColorScheme / etc. all
# resolve to that same theme part), so we rewrite its twelve colour
# slots in place to the standard Office theme values.

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
